# Weekly data refresh: insert one new price record for "Poroto verde"
# (Vega Central Mapocho de Santiago) at row 230 of Sheet1, pushing the
# existing rows 230-248 down to 231-249 (dimension grows from R248 to R249).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 230, shifting rows 230:248 -> 231:249.
$ws.Rows.Item(230).Insert()

# Populate the newly inserted row 230 with the new weekly record.
$ws.Range("A230").Value = 9
$ws.Range("B230").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C230").Value = "Metropolitana"
$ws.Range("D230").Value = 44461
$ws.Range("E230").Value = 13
$ws.Range("F230").Value = 100112031
$ws.Range("G230").Value = "Poroto verde"
$ws.Range("H230").Value = "Magnum"
$ws.Range("I230").Value = "Primera"
$ws.Range("J230").Value = 25
$ws.Range("K230").Value = 40000
$ws.Range("L230").Value = 42000
$ws.Range("M230").Value = 41040
$ws.Range("N230").Value = "`$/malla 25 kilos"
$ws.Range("O230").Value = "Perú"
$ws.Range("P230").Value = 1642
$ws.Range("Q230").Value = 25
$ws.Range("R230").Value = "Hortaliza"
